$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.288.41'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').Value = '1.689.95'
$ws.Range('E3').Value = '  +1.29%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.006'
$ws.Range('D4').Style = "Normal"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '218.87'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.77%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2697'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +1.76%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06441'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.57%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '22.05'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +2.11%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07471'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.48%  '
$ws.Range('D12').Value = '1.695.67'
$ws.Range('E12').Value = '  +1.52%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.559'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.34%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.5858'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +1.43%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.000008528'
$ws.Range('D15').Style = "Normal"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '64.64'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.49%  '
$ws.Range('D17').Value = '26.322.34'
$ws.Range('E17').Value = '  +0.30%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '4.962'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.69%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.006'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '10.87'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.52%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '189.69'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.43%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.221'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.64%  '
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '144.76'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.16%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '7.673'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.26%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.1233'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +5.23%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '15.85'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.96%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.06702'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +15.96%  '
$ws.Range('E29').Value = '  +5.61%  '
$ws.Range('E30').Value = '  +1.05%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.589'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +2.42%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.567'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +1.31%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.670'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.42%  '
$ws.Range('E34').Value = '  +2.40%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.6220'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +3.98%  '
$ws.Range('E36').Value = '  +1.57%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.699'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +2.30%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '6.329'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +5.47%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01623'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.92%  '
$ws.Range('D40').Value = '1.105.29'
$ws.Range('E40').Value = '  +2.27%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.8851'
$ws.Range('D41').Style = "Normal"
$ws.Range('E42').Value = '  +0.85%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '101.31'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +1.45%  '
$ws.Range('D44').Value = '1.836.96'
$ws.Range('E44').Value = '  +1.17%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.00000000112'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.84%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '56.85'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.99%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '8.166'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +1.52%  '
$ws.Range('E48').Value = '  -0.08%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.05264'
$ws.Range('D49').Style = "Normal"
$ws.Range('E50').Value = '  +0.06%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '6.066'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +3.57%  '
